$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Body text edit: the diary entry "...创建了一个dev分支。" gets a new
#    sentence appended right after it: "使用Git创建分支简单又便捷。"
#    Word drops its hidden "_GoBack" bookmark at the location of the most
#    recent edit, so it moves from the empty paragraph below into this
#    paragraph, landing right before the freshly typed trailing "。".
# ---------------------------------------------------------------------------

$wdFindContinue = 1
$wdReplaceNone  = 0
$wdCollapseEnd  = 0
$wdCharacter    = 1

# Locate the sentence by its distinctive text instead of hard-coded
# character offsets, so the edit is resilient to anything earlier in the doc.
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("创建了一个dev分支。", $true, $false, $false, $false, $false, $true, `
                               $wdFindContinue, $false, "", $wdReplaceNone)

if ($found) {
    # Collapse to right after the matched sentence (after its "。") and type
    # the new sentence there.
    [void]$anchor.Collapse($wdCollapseEnd)

    $newSentence = "使用Git创建分支简单又便捷。"
    $anchor.InsertAfter($newSentence)

    # $anchor now sits collapsed right after the inserted text. Re-home the
    # _GoBack bookmark just before the last typed character (the new
    # trailing "。"), matching where Word leaves it after typing.
    [void]$anchor.MoveEnd($wdCharacter, -1)
    [void]$anchor.Collapse($wdCollapseEnd)

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $gobackRange = $d.Range($anchor.End, $anchor.End)
    $d.Bookmarks.Add("_GoBack", $gobackRange)
}

# ---------------------------------------------------------------------------
# 2) Style metadata tweak: mark the built-in "Default Paragraph Font"
#    character style as a Quick Style (adds <w:qFormat/> to its definition).
# ---------------------------------------------------------------------------
$defaultParaFont = $d.Styles.Item("Default Paragraph Font")
$defaultParaFont.QuickStyle = $true
